$p = $ppt.ActivePresentation

# 1) Add a new blank slide at the end (slide 18), using the "Blank" layout (ppLayoutBlank = 12)
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 12)

# 2) Merge the three runs of the last paragraph on slide 9 ("Existe " + "desconcentración " +
#    "por parte las redes sociales, debido a la conexión a internet") into a single run.
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(6, 1)

$run1 = $para.Runs(1, 1)
$run2 = $para.Runs(2, 1)
$run3 = $para.Runs(3, 1)
$fullText = $run1.Text + $run2.Text + $run3.Text

$start2 = $run2.Start
$len23 = $run2.Text.Length + $run3.Text.Length
$tail = $tr.Characters($start2, $len23)
$tail.Delete()

$mergedRun = $tr.Paragraphs(6, 1).Runs(1, 1)
$mergedRun.Text = $fullText
